$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '52.124.75'
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$ws.Range("D3").Value = '2.995.75'
$ws.Range("E3").Value = '  +2.60%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.24%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.22'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.63%  '

# Row 11
$ws.Range("E11").Value = '  +1.77%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0855'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.67%  '

# Row 14
$ws.Range("D14").Value = '3.466.20'
$ws.Range("E14").Value = '  +2.57%  '

# Row 15
$ws.Range("E15").Value = '  -4.29%  '

# Row 16
$ws.Range("D16").Value = '2.998.95'
$ws.Range("E16").Value = '  +2.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.81%  '

# Row 18
$ws.Range("D18").Value = '52.125.16'
$ws.Range("E18").Value = '  +0.92%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.76%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.10%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0971'
$ws.Range("E22").Value = '  -0.83%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.19%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.176'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.96%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.37%  '

# Row 29
$ws.Range("E29").Value = '  -0.01%  '

# Row 30
$ws.Range("E30").Value = '  -1.94%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.38%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.29%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '36.07'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.67%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.25%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '51.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.21%  '

# Row 36
$ws.Range("E36").Value = '  +0.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.78%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.52%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.47%  '

# Row 41
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.20%  '

# Row 42
$ws.Range("E42").Value = '  -0.36%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '124.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.70%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.13%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.56%  '

# Row 46
$ws.Range("D46").Value = '2.121.88'
$ws.Range("E46").Value = '  -0.74%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.13%  '

# Row 48
$ws.Range("E48").Value = '  -5.97%  '

# Row 49
$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.243'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.19%  '

# Row 50
$ws.Range("B50").Value = 'BEAM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0331'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.67%  '

# Row 51
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.906'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.51%  '
